$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the subheader row (row 2: ACCESSION NO / TITLE / ACQUISITION NO. / Item /
# ITEM DESCRIPTION / LOCATION | SECTION / QTY / AR NUMBER / 2022/23 RFID Number / COLLECTIONS).
# This shifts every data row up by one.
$ws.Rows(2).Delete()

# The QTY / AR NUMBER / 2022/23 RFID Number / COLLECTIONS headers that used to live in the
# (now deleted) row 2 for the helper columns K:N need to move up into row 1, matching the
# other column headers.
$ws.Range("K1").Value = "QTY"
$ws.Range("L1").Value = "AR NUMBER"
$ws.Range("M1").Value = "2022/23 RFID Number"
$ws.Range("N1").Value = "COLLECTIONS"
$ws.Range("K1:N1").Font.Bold = $true

# Remove the acquisition date ("ACQUISITION NO.") values from column D - the collection no
# longer tracks an acquisition date.
$ws.Range("D2:D6").ClearContents()

# Un-hide the helper columns now that they carry a visible header in row 1.
$ws.Columns("K:N").Hidden = $false

# Refresh the column widths to better fit the remaining visible content.
$ws.Columns("G:I").AutoFit()

[void]$ws.Range("D2:D6").Select()

$wb.Save()
